# changement horaire et responsabilites
# Insert a new row for "Gabriel Montplaisir" into the sorted team roster,
# placed alphabetically between "Frédérik Taleb" (row 8) and
# "Mathieu Fréchette" (the old row 9, which becomes row 10).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a fresh row at position 9; this shifts the former rows 9-11
#    (Mathieu Fréchette, Sébastien Trottier, Simon Tousignant) down to 10-12.
$ws.Rows("9:9").Insert()

# 2) Populate the new row 9 with Gabriel Montplaisir's info (values first,
#    formats copied in afterwards - doing it in the other order strips the
#    quote-prefixed "Hyperlink" style back down to its plain variant).
$ws.Range("A9").Value = "Gabriel Montplaisir"
$ws.Range("B9").Value = ""
$ws.Range("C9").Value = "/assets/placeholder-man.webp"
$ws.Range("D9").Value = "montplaisir.gabriel@cegepvicto.ca"
$ws.Range("E9").Value = "(819) 758-6401 poste 2519"
$ws.Range("F9").Value = "C-207"
$ws.Range("G9").Value = "https://teams.microsoft.com/l/chat/0/0?tenantId=9d6cf526-ad81-46f8-a73a-a507aaf06cda&users=MONTPLAISIR.GABRIEL@cegepvicto.ca"

# 3) Match the cell styles used by the analogous "Simon Tousignant" row
#    (C -> avatar/hyperlink style, D -> quoted hyperlink style), which is
#    the template this new row was copied from.
$ws.Range("C12").Copy()
$ws.Range("C9").PasteSpecial(-4122)
$ws.Range("D12").Copy()
$ws.Range("D9").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# 4) Rebuild every hyperlink in the sheet: inserting the row shifted the
#    underlying data down but this environment does not auto-shift the
#    hyperlink anchors, so clear them all and re-create with the same
#    targets pointing at their (possibly new) cells.
$ws.Range("D6").Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("D6"), "mailto:rivard.etienne@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G6"), "https://teams.microsoft.com/l/chat/0/0?tenantId=9d6cf526-ad81-46f8-a73a-a507aaf06cda&users=RIVARD.ETIENNE@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:ouellet.alexandre@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G2"), "https://teams.microsoft.com/l/chat/0/0?tenantId=9d6cf526-ad81-46f8-a73a-a507aaf06cda&users=OUELLET.ALEXANDRE@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D3"), "mailto:croteau.carine@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G3"), "https://teams.microsoft.com/l/chat/0/0?tenantId=9d6cf526-ad81-46f8-a73a-a507aaf06cda&users=CROTEAU.CARINE@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D10"), "mailto:frechette.mathieu@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G10"), "https://teams.microsoft.com/l/chat/0/0?tenantId=9d6cf526-ad81-46f8-a73a-a507aaf06cda&users=FRECHETTE.MATHIEU@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "mailto:lagace.christiane@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G4"), "https://teams.microsoft.com/l/chat/0/0?tenantId=9d6cf526-ad81-46f8-a73a-a507aaf06cda&users=LAGACE.CHRISTIANE@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D7"), "mailto:mercier.francois@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G7"), "https://teams.microsoft.com/l/chat/0/0?tenantId=9d6cf526-ad81-46f8-a73a-a507aaf06cda&users=MERCIER.FRANCOIS@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D8"), "mailto:taleb.frederik@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G8"), "https://teams.microsoft.com/l/chat/0/0?tenantId=9d6cf526-ad81-46f8-a73a-a507aaf06cda&users=TALEB.FREDERIK@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D11"), "mailto:tousignant.simon@cegepvicto.ca", "", "", "tousignant.simon@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G11"), "https://teams.microsoft.com/l/chat/0/0?tenantId=9d6cf526-ad81-46f8-a73a-a507aaf06cda&users=TROTTIER.SEBASTIEN@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "mailto:chaieb.cirine@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G5"), "https://teams.microsoft.com/l/chat/0/0?tenantId=9d6cf526-ad81-46f8-a73a-a507aaf06cda&users=CHAIEB.CIRINE@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D12"), "mailto:tousignant.simon@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G12"), "https://teams.microsoft.com/l/chat/0/0?tenantId=9d6cf526-ad81-46f8-a73a-a507aaf06cda&users=TOUSIGNANT.SIMON@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C12"), "https://avatar.iran.liara.run/public/boy", "", "", "https://avatar.iran.liara.run/public/boy") | Out-Null

# New hyperlinks for Gabriel Montplaisir's row (copied from the Simon
# Tousignant template row, hence the stale display text on D9).
$ws.Hyperlinks.Add($ws.Range("D9"), "mailto:montplaisir.gabriel@cegepvicto.ca", "", "", "tousignant.simon@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G9"), "https://teams.microsoft.com/l/chat/0/0?tenantId=9d6cf526-ad81-46f8-a73a-a507aaf06cda&users=MONTPLAISIR.GABRIEL@cegepvicto.ca") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C9"), "https://avatar.iran.liara.run/public/boy", "", "", "https://avatar.iran.liara.run/public/boy") | Out-Null

# 5) Cosmetic bits that followed from the edit: selection moved to G10.
$ws.Range("G10").Select()

"done"
